$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in the title cell (A1).
$ws.Range("A1").Value = "Datos actualizados a 26 de Julio de 2020 a las 02:58"

# Update country case-count rows with the latest data. Updating the country
# name cell as well since several countries swap position in the (descending,
# by total cases) ranking, changing which country name lands on which row.

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 4315683
$ws.Cells.Item(4, 3).Value = 67387
$ws.Cells.Item(4, 4).Value = 2061692
$ws.Cells.Item(4, 5).Value = 2104596
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 905
$ws.Cells.Item(4, 8).Value = 149395

# Row 24: Canada
$ws.Cells.Item(24, 1).Value = "Canada"
$ws.Cells.Item(24, 2).Value = 113556
$ws.Cells.Item(24, 3).Value = 350
$ws.Cells.Item(24, 4).Value = 99115
$ws.Cells.Item(24, 5).Value = 5556
$ws.Cells.Item(24, 6).Value = 0
$ws.Cells.Item(24, 7).Value = 4
$ws.Cells.Item(24, 8).Value = 8885

# Row 42: Panama
$ws.Cells.Item(42, 1).Value = "Panama"
$ws.Cells.Item(42, 2).Value = 58864
$ws.Cells.Item(42, 3).Value = 871
$ws.Cells.Item(42, 4).Value = 33428
$ws.Cells.Item(42, 5).Value = 24161
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(42, 7).Value = 25
$ws.Cells.Item(42, 8).Value = 1275

# Row 43: Emiratos Arabes Unidos
$ws.Cells.Item(43, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(43, 2).Value = 58562
$ws.Cells.Item(43, 3).Value = 313
$ws.Cells.Item(43, 4).Value = 51628
$ws.Cells.Item(43, 5).Value = 6591
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(43, 7).Value = 0
$ws.Cells.Item(43, 8).Value = 343

# Row 119: Cuba
$ws.Cells.Item(119, 1).Value = "Cuba"
$ws.Cells.Item(119, 2).Value = 2478
$ws.Cells.Item(119, 3).Value = 9
$ws.Cells.Item(119, 4).Value = 2345
$ws.Cells.Item(119, 5).Value = 46
$ws.Cells.Item(119, 6).Value = 0
$ws.Cells.Item(119, 7).Value = 0
$ws.Cells.Item(119, 8).Value = 87

# Row 122: Cabo Verde
$ws.Cells.Item(122, 1).Value = "Cabo Verde"
$ws.Cells.Item(122, 2).Value = 2258
$ws.Cells.Item(122, 3).Value = 38
$ws.Cells.Item(122, 4).Value = 1363
$ws.Cells.Item(122, 5).Value = 873
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 7).Value = 0
$ws.Cells.Item(122, 8).Value = 22

# Row 146: Burkina Faso
$ws.Cells.Item(146, 1).Value = "Burkina Faso"
$ws.Cells.Item(146, 2).Value = 1086
$ws.Cells.Item(146, 3).Value = 11
$ws.Cells.Item(146, 4).Value = 920
$ws.Cells.Item(146, 5).Value = 113
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 53

# Row 147: Republica de Chipre
$ws.Cells.Item(147, 1).Value = "Republica de Chipre"
$ws.Cells.Item(147, 2).Value = 1053
$ws.Cells.Item(147, 3).Value = 6
$ws.Cells.Item(147, 4).Value = 852
$ws.Cells.Item(147, 5).Value = 182
$ws.Cells.Item(147, 6).Value = 0
$ws.Cells.Item(147, 7).Value = 0
$ws.Cells.Item(147, 8).Value = 19

# Row 151: Santo Tome y Principe
$ws.Cells.Item(151, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(151, 2).Value = 862
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(151, 4).Value = 662
$ws.Cells.Item(151, 5).Value = 186
$ws.Cells.Item(151, 6).Value = 0
$ws.Cells.Item(151, 7).Value = 0
$ws.Cells.Item(151, 8).Value = 14

# Row 164: Burundi
$ws.Cells.Item(164, 1).Value = "Burundi"
$ws.Cells.Item(164, 2).Value = 361
$ws.Cells.Item(164, 3).Value = 16
$ws.Cells.Item(164, 4).Value = 279
$ws.Cells.Item(164, 5).Value = 81
$ws.Cells.Item(164, 6).Value = 0
$ws.Cells.Item(164, 7).Value = 0
$ws.Cells.Item(164, 8).Value = 1

# Row 165: Guyana
$ws.Cells.Item(165, 1).Value = "Guyana"
$ws.Cells.Item(165, 2).Value = 360
$ws.Cells.Item(165, 3).Value = 8
$ws.Cells.Item(165, 4).Value = 180
$ws.Cells.Item(165, 5).Value = 160
$ws.Cells.Item(165, 6).Value = 0
$ws.Cells.Item(165, 7).Value = 0
$ws.Cells.Item(165, 8).Value = 20

# Row 166: Birmania
$ws.Cells.Item(166, 1).Value = "Birmania"
$ws.Cells.Item(166, 2).Value = 348
$ws.Cells.Item(166, 3).Value = 2
$ws.Cells.Item(166, 4).Value = 288
$ws.Cells.Item(166, 5).Value = 54
$ws.Cells.Item(166, 6).Value = 0
$ws.Cells.Item(166, 7).Value = 0
$ws.Cells.Item(166, 8).Value = 6

# Row 170: Bahamas
$ws.Cells.Item(170, 1).Value = "Bahamas"
$ws.Cells.Item(170, 2).Value = 326
$ws.Cells.Item(170, 3).Value = 10
$ws.Cells.Item(170, 4).Value = 91
$ws.Cells.Item(170, 5).Value = 224
$ws.Cells.Item(170, 6).Value = 0
$ws.Cells.Item(170, 7).Value = 0
$ws.Cells.Item(170, 8).Value = 11

# Row 172: Gambia
$ws.Cells.Item(172, 1).Value = "Gambia"
$ws.Cells.Item(172, 2).Value = 277
$ws.Cells.Item(172, 3).Value = 61
$ws.Cells.Item(172, 4).Value = 60
$ws.Cells.Item(172, 5).Value = 211
$ws.Cells.Item(172, 6).Value = 0
$ws.Cells.Item(172, 7).Value = 0
$ws.Cells.Item(172, 8).Value = 6

# Row 173: Martinica
$ws.Cells.Item(173, 1).Value = "Martinica"
$ws.Cells.Item(173, 2).Value = 269
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(173, 4).Value = 98
$ws.Cells.Item(173, 5).Value = 156
$ws.Cells.Item(173, 6).Value = 0
$ws.Cells.Item(173, 7).Value = 0
$ws.Cells.Item(173, 8).Value = 15

# Row 174: Eritrea
$ws.Cells.Item(174, 1).Value = "Eritrea"
$ws.Cells.Item(174, 2).Value = 263
$ws.Cells.Item(174, 3).Value = 2
$ws.Cells.Item(174, 4).Value = 189
$ws.Cells.Item(174, 5).Value = 74
$ws.Cells.Item(174, 6).Value = 0
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 0

# Row 175: Camboya
$ws.Cells.Item(175, 1).Value = "Camboya"
$ws.Cells.Item(175, 2).Value = 225
$ws.Cells.Item(175, 3).Value = 23
$ws.Cells.Item(175, 4).Value = 143
$ws.Cells.Item(175, 5).Value = 82
$ws.Cells.Item(175, 6).Value = 0
$ws.Cells.Item(175, 7).Value = 0
$ws.Cells.Item(175, 8).Value = 0

# Row 185: Seychelles
$ws.Cells.Item(185, 1).Value = "Seychelles"
$ws.Cells.Item(185, 2).Value = 114
$ws.Cells.Item(185, 3).Value = 0
$ws.Cells.Item(185, 4).Value = 39
$ws.Cells.Item(185, 5).Value = 75
$ws.Cells.Item(185, 6).Value = 0
$ws.Cells.Item(185, 7).Value = 0
$ws.Cells.Item(185, 8).Value = 0

# Row 210: Islas Malvinas
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 2).Value = 13
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 13
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Groenlandia
$ws.Cells.Item(211, 1).Value = "Groenlandia"
$ws.Cells.Item(211, 2).Value = 13
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 13
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 0
